$d = $word.ActiveDocument

$pairs = @(
    @("97÷9=", "15÷8="),
    @("19÷6=", "77÷9="),
    @("46÷9=", "97÷8="),
    @("35÷5=", "73÷9="),
    @("75÷9=", "30÷5="),
    @("61÷5=", "86÷8="),
    @("25÷3=", "44÷8="),
    @("75÷8=", "96÷2="),
    @("91÷6=", "20÷3="),
    @("23÷2=", "54÷6="),
    @("82÷3=", "66÷7="),
    @("68÷2=", "95÷6="),
    @("96÷5=", "64÷6="),
    @("36÷8=", "49÷5="),
    @("45÷8=", "13÷7="),
    @("64÷8=", "10÷7="),
    @("49÷7=", "95÷6="),
    @("21÷3=", "87÷9="),
    @("64÷2=", "13÷6="),
    @("13÷8=", "87÷4="),
    @("43÷4=", "59÷4="),
    @("50÷2=", "32÷5="),
    @("87÷5=", "98÷6="),
    @("94÷9=", "38÷3="),
    @("77÷6=", "34÷4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
